# Restore / update the vendor schedule table (rows re-sorted by vendor
# name and several status/date values refreshed), rename the
# "Rockwell" vendor entry to "Gigaset", and mark the two note cells
# (C9/C10) as wrap-text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Siemens ---------------------------------------------------
$ws.Range("A2").Value = "Siemens"
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 44894
$ws.Range("D2").Value = 44993
$ws.Range("E2").ClearContents()

# --- Row 3: Asus --------------------------------------------------------
$ws.Range("A3").Value = "Asus"
$ws.Range("B3").Value = 100
$ws.Range("C3").Value = 44893
$ws.Range("D3").Value = 45261

# --- Row 4: Schneider ----------------------------------------------------
$ws.Range("A4").Value = "Schneider"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 44901
$ws.Range("D4").Value = 45261
$ws.Range("E4").Value = "SchneiderElectricScraper"
$ws.Range("F4").Value = 10

# --- Row 5: AVM ----------------------------------------------------------
$ws.Range("A5").Value = "AVM"
$ws.Range("B5").Value = 100
$ws.Range("C5").Value = 44902
$ws.Range("D5").Value = 45261
$ws.Range("E5").Value = "AVMScraper"
$ws.Range("F5").ClearContents()

# --- Row 6: Synology -------------------------------------------------------
$ws.Range("A6").Value = "Synology"
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = 44902
$ws.Range("D6").Value = 45261
$ws.Range("E6").Value = "SynologyScraper"

# --- Row 7: Gigaset (was Rockwell) -----------------------------------------
$ws.Range("A7").Value = "Gigaset"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 44934
$ws.Range("D7").Value = 44943
$ws.Range("E7").Value = "GigasetScraper"

# --- Row 8: Swisscom ---------------------------------------------------
$ws.Range("A8").Value = "Swisscom"
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = 44902
$ws.Range("D8").Value = 44927
$ws.Range("E8").Value = "SwisscomScraper"

# --- Row 9: Zyxel --------------------------------------------------------
$ws.Range("A9").Value = "Zyxel"
$ws.Range("B9").Value = 100
$ws.Range("C9").WrapText = $true
$ws.Range("D9").Value = 44927
$ws.Range("E9").Value = "ZyxelScraper"

# --- Row 10: ABB ---------------------------------------------------------
$ws.Range("A10").Value = "ABB"
$ws.Range("B10").Value = 100
$ws.Range("C10").WrapText = $true
$ws.Range("D10").Value = 44927
$ws.Range("E10").Value = "ABBScraper"

# --- Row 11: Trendnet ------------------------------------------------------
$ws.Range("A11").Value = "Trendnet"
$ws.Range("B11").Value = 100
$ws.Range("C11").Value = 44902
$ws.Range("D11").Value = 44927
$ws.Range("E11").Value = "TrendnetScraper"

# --- Row 12: TP-Link (unchanged) -------------------------------------------
# left as-is

# --- Sheet selection -----------------------------------------------------
$ws.Range("I22").Select() | Out-Null
